$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Person" sheet: drop the Relationship (E) and Education (I) columns,
#    insert a new Hobbies column, split Complexion/HairKind into two
#    columns and extend the HairKind list.
# ------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")

# Remove the (now redundant) Education column first so the later
# deletion of Relationship doesn't shift its address.
$wsPerson.Columns("I:I").Delete()
$wsPerson.Columns("E:E").Delete()

# Make room for the new "Hobbies" column between LanguagesKnown and
# Complexion.
$wsPerson.Columns("I:I").Insert()

# Re-label the (renamed) HairKind column and extend its value list -
# order matters here: it controls shared-string allocation order.
$wsPerson.Range("K2").Value = "HairKind"
$wsPerson.Range("K3").Value = "Thick"
$wsPerson.Range("K4").Value = "Sparse"
$wsPerson.Range("K5").Value = "BaldHead"
$wsPerson.Range("I2").Value = "Hobbies"

$wsPerson.Range("L5").Select()

# ------------------------------------------------------------------
# 2) "Family" sheet: nothing content-wise changes, but touch the view
#    so the stale topLeftCell scroll position is cleared.
# ------------------------------------------------------------------
$wsFamily = $wb.Worksheets.Item("Family")
$wsFamily.Range("S1").Select()

# ------------------------------------------------------------------
# 3) New "Location" sheet, appended after "Family".
# ------------------------------------------------------------------
$wsLocation = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsLocation.Name = "Location"
$wsLocation.Range("A1").Value = "loginId"
$wsLocation.Range("B1").Value = "Country"
$wsLocation.Range("C1").Value = "State"
$wsLocation.Range("D1").Value = "City"
$wsLocation.Range("E1").Value = "Address"
$wsLocation.Range("F1").Value = "PostalCode"
$wsLocation.Range("A2").Select()

# ------------------------------------------------------------------
# 4) "Education" sheet becomes the active tab, with B3:B22 selected.
# ------------------------------------------------------------------
$wsEducation = $wb.Worksheets.Item("Education")
$wsEducation.Range("B3:B22").Select()
$wsEducation.Activate()
